# Abundance estimate data update: append Week 22 rows to Sheet1 (per-stratum
# counts) and Sheet2 (week/date lookup), matching the weekly refresh pattern
# used for every prior week in this tracking workbook.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# ---------------------------------------------------------------------
# Sheet1: seven new rows (144-150) for Week 22 - one row per stratum plus
# the "All Strata" summary row, same column layout as every prior week.
# ---------------------------------------------------------------------
$week22 = @(
    @{ Row = 144; Stratum = "Suisun Bay";          Sites = 6;  Tows = 24;  Volume = 93861 },
    @{ Row = 145; Stratum = "Suisun Marsh";         Sites = 6;  Tows = 24;  Volume = 95286 },
    @{ Row = 146; Stratum = "Lower Sacramento";     Sites = 6;  Tows = 24;  Volume = 94751 },
    @{ Row = 147; Stratum = "Cache Slough LI";      Sites = 6;  Tows = 24;  Volume = 80263 },
    @{ Row = 148; Stratum = "Sac DW Ship Channel";  Sites = 6;  Tows = 24;  Volume = 80985 },
    @{ Row = 149; Stratum = "Lower San Joaquin";    Sites = 6;  Tows = 24;  Volume = 92989 },
    @{ Row = 150; Stratum = "All Strata";           Sites = 36; Tows = 144; Volume = 538136 }
)

foreach ($rowData in $week22) {
    $r = $rowData.Row
    $ws1.Range("A$r").Value = 22
    $ws1.Range("B$r").Value = $rowData.Stratum
    $ws1.Range("C$r").Value = $rowData.Sites
    $ws1.Range("D$r").Value = $rowData.Tows
    $ws1.Range("E$r").Value = 0
    $ws1.Range("F$r").Value = 0
    $ws1.Range("G$r").Value = 0
    $ws1.Range("H$r").Value = $rowData.Volume
    $ws1.Range("H$r").NumberFormat = "#,##0"
    $ws1.Range("I$r").Value = "0*"
    $ws1.Range("J$r").Value = "NA"
    $ws1.Range("K$r").Value = "NA"
}

# ---------------------------------------------------------------------
# Sheet2: new week/date lookup row for Week 22.
# ---------------------------------------------------------------------
$ws2.Range("A23").Value = 22
$ws2.Range("B23").Value = "October 27–31, 2025"

# ---------------------------------------------------------------------
# Restore the on-screen selections that come along with this edit: the
# last-touched cell on Sheet1 and the active cell on Sheet2, leaving
# Sheet2 as the active tab (matches the saved workbook state).
# ---------------------------------------------------------------------
$ws1.Range("L140").Select() | Out-Null
$ws2.Range("I14").Select() | Out-Null
